$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1081.579
$ws.Cells.Item(40, 9).Value = 971.875
$ws.Cells.Item(40, 10).Value = 1666.6666
$ws.Cells.Item(40, 11).Value = 971.875
$ws.Cells.Item(40, 12).Value = 1666.6666
$ws.Cells.Item(40, 13).Value = -796.875
$ws.Cells.Item(40, 14).Value = -2016.6666
$ws.Cells.Item(135, 8).Value = 48183172
$ws.Cells.Item(135, 9).Value = 758.4
$ws.Cells.Item(135, 10).Value = 88335180
$ws.Cells.Item(135, 11).Value = 6825.599999999999
$ws.Cells.Item(135, 12).Value = 795016620
$ws.Cells.Item(135, 13).Value = -4290.599999999999
$ws.Cells.Item(135, 14).Value = -795021690
$ws.Cells.Item(138, 8).Value = 2459.15
$ws.Cells.Item(138, 9).Value = 1648.25
$ws.Cells.Item(138, 10).Value = 2492.9375
$ws.Cells.Item(138, 11).Value = 4944.75
$ws.Cells.Item(138, 12).Value = 7478.8125
$ws.Cells.Item(138, 13).Value = 195.25
$ws.Cells.Item(138, 14).Value = -17758.8125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9645.52
$ws.Cells.Item(32, 9).Value = 8169.9316
$ws.Cells.Item(32, 11).Value = 8169.9316
$ws.Cells.Item(32, 13).Value = -7882.9316
$ws.Cells.Item(74, 8).Value = 20005578
$ws.Cells.Item(74, 9).Value = 31579770
$ws.Cells.Item(74, 10).Value = 13790.728
$ws.Cells.Item(74, 11).Value = 31579770
$ws.Cells.Item(74, 12).Value = 13790.728
$ws.Cells.Item(74, 13).Value = -31578896
$ws.Cells.Item(74, 14).Value = -15538.728
$ws.Cells.Item(77, 8).Value = 20005578
$ws.Cells.Item(77, 9).Value = 31579770
$ws.Cells.Item(77, 10).Value = 13790.728
$ws.Cells.Item(77, 11).Value = 157898850
$ws.Cells.Item(77, 12).Value = 68953.64
$ws.Cells.Item(77, 13).Value = -157894482
$ws.Cells.Item(77, 14).Value = -77689.64
$ws.Cells.Item(132, 8).Value = 988509.9399999999
$ws.Cells.Item(132, 9).Value = 1234463.6
$ws.Cells.Item(132, 10).Value = 127672.25
$ws.Cells.Item(132, 11).Value = 3703390.8
$ws.Cells.Item(132, 12).Value = 383016.75
$ws.Cells.Item(132, 13).Value = -3700860.8
$ws.Cells.Item(132, 14).Value = -388076.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2156.1396
$ws.Cells.Item(86, 9).Value = 1797.7222
$ws.Cells.Item(86, 10).Value = 3999.4285
$ws.Cells.Item(86, 11).Value = 1797.7222
$ws.Cells.Item(86, 12).Value = 3999.4285
$ws.Cells.Item(86, 13).Value = -674.7221999999999
$ws.Cells.Item(86, 14).Value = -6245.4285
$ws.Cells.Item(89, 8).Value = 2156.1396
$ws.Cells.Item(89, 9).Value = 1797.7222
$ws.Cells.Item(89, 10).Value = 3999.4285
$ws.Cells.Item(89, 11).Value = 8988.610999999999
$ws.Cells.Item(89, 12).Value = 19997.1425
$ws.Cells.Item(89, 13).Value = -3372.610999999999
$ws.Cells.Item(89, 14).Value = -31229.1425

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(141, 8).Value = 10740.556
$ws.Cells.Item(141, 9).Value = 8333
$ws.Cells.Item(141, 10).Value = 11428.429
$ws.Cells.Item(141, 11).Value = 8333
$ws.Cells.Item(141, 12).Value = 11428.429
$ws.Cells.Item(141, 13).Value = -3153
$ws.Cells.Item(141, 14).Value = -21788.429

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 50129.9
$ws.Cells.Item(12, 9).Value = 76.59999999999999
$ws.Cells.Item(12, 10).Value = 100183.2
$ws.Cells.Item(12, 11).Value = 229.8
$ws.Cells.Item(12, 12).Value = 300549.6
$ws.Cells.Item(12, 13).Value = -56.79999999999998
$ws.Cells.Item(12, 14).Value = -300895.6
$ws.Cells.Item(17, 8).Value = 1270
$ws.Cells.Item(17, 9).Value = 540
$ws.Cells.Item(17, 10).Value = 2000
$ws.Cells.Item(17, 11).Value = 1620
$ws.Cells.Item(17, 12).Value = 6000
$ws.Cells.Item(17, 13).Value = -1451
$ws.Cells.Item(17, 14).Value = -6338
$ws.Cells.Item(34, 8).Value = 689.375
$ws.Cells.Item(34, 9).Value = 267.14285
$ws.Cells.Item(34, 10).Value = 1017.7778
$ws.Cells.Item(34, 11).Value = 801.4285500000001
$ws.Cells.Item(34, 12).Value = 3053.3334
$ws.Cells.Item(34, 13).Value = -717.4285500000001
$ws.Cells.Item(34, 14).Value = -3221.3334
$ws.Cells.Item(39, 8).Value = 12444.444
$ws.Cells.Item(39, 10).Value = 12444.444
$ws.Cells.Item(39, 12).Value = 37333.33199999999
$ws.Cells.Item(39, 14).Value = -37921.33199999999
$ws.Cells.Item(55, 8).Value = 2563.6365
$ws.Cells.Item(55, 9).Value = 100
$ws.Cells.Item(55, 10).Value = 2810
$ws.Cells.Item(55, 11).Value = 300
$ws.Cells.Item(55, 12).Value = 8430
$ws.Cells.Item(55, 13).Value = -123
$ws.Cells.Item(55, 14).Value = -8784
$ws.Cells.Item(113, 8).Value = 2481.14
$ws.Cells.Item(113, 9).Value = 503.94736
$ws.Cells.Item(113, 10).Value = 3692.9678
$ws.Cells.Item(113, 11).Value = 1511.84208
$ws.Cells.Item(113, 12).Value = 11078.9034
$ws.Cells.Item(113, 13).Value = 658.1579200000001
$ws.Cells.Item(113, 14).Value = -15418.9034
$ws.Cells.Item(126, 8).Value = 3016.5
$ws.Cells.Item(126, 9).Value = 1006.8571
$ws.Cells.Item(126, 10).Value = 5830
$ws.Cells.Item(126, 11).Value = 3020.5713
$ws.Cells.Item(126, 12).Value = 17490
$ws.Cells.Item(126, 13).Value = 1919.4287
$ws.Cells.Item(126, 14).Value = -27370

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(18, 8).Value = 3111.111
$ws.Cells.Item(18, 10).Value = 4000
$ws.Cells.Item(18, 12).Value = 4000
$ws.Cells.Item(18, 14).Value = -4586
$ws.Cells.Item(43, 8).Value = 2758.3333
$ws.Cells.Item(43, 9).Value = 1012.5
$ws.Cells.Item(43, 10).Value = 6250
$ws.Cells.Item(43, 11).Value = 1012.5
$ws.Cells.Item(43, 12).Value = 6250
$ws.Cells.Item(43, 13).Value = -861.5
$ws.Cells.Item(43, 14).Value = -6552
$ws.Cells.Item(46, 8).Value = 24750
$ws.Cells.Item(46, 10).Value = 24750
$ws.Cells.Item(46, 12).Value = 24750
$ws.Cells.Item(46, 14).Value = -25062
$ws.Cells.Item(57, 8).Value = 12487.5
$ws.Cells.Item(57, 10).Value = 12487.5
$ws.Cells.Item(57, 12).Value = 12487.5
$ws.Cells.Item(57, 14).Value = -14127.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 2536.348
$ws.Cells.Item(136, 9).Value = 1580.375
$ws.Cells.Item(136, 10).Value = 4721.4287
$ws.Cells.Item(136, 11).Value = 4741.125
$ws.Cells.Item(136, 12).Value = 14164.2861
$ws.Cells.Item(136, 13).Value = -2191.125
$ws.Cells.Item(136, 14).Value = -19264.2861

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 5365.6665
$ws.Cells.Item(81, 9).Value = 2587.75
$ws.Cells.Item(81, 10).Value = 7588
$ws.Cells.Item(81, 11).Value = 5175.5
$ws.Cells.Item(81, 12).Value = 15176
$ws.Cells.Item(81, 13).Value = -4114.5
$ws.Cells.Item(81, 14).Value = -17298
$ws.Cells.Item(84, 8).Value = 5365.6665
$ws.Cells.Item(84, 9).Value = 2587.75
$ws.Cells.Item(84, 10).Value = 7588
$ws.Cells.Item(84, 11).Value = 25877.5
$ws.Cells.Item(84, 12).Value = 75880
$ws.Cells.Item(84, 13).Value = -20573.5
$ws.Cells.Item(84, 14).Value = -86488
$ws.Cells.Item(136, 8).Value = 6469.8335
$ws.Cells.Item(136, 9).Value = 8052.185
$ws.Cells.Item(136, 10).Value = 1722.7778
$ws.Cells.Item(136, 11).Value = 24156.555
$ws.Cells.Item(136, 12).Value = 5168.3334
$ws.Cells.Item(136, 13).Value = -21606.555
$ws.Cells.Item(136, 14).Value = -10268.3334
